# Insert a new row at position 140, pushing existing rows 140-167 down to 141-168.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("140:140").Insert()
# The Insert() operation copies formatting from the row above; the new row
# should be unstyled (matching the target XML, which has no "s" attribute).
$ws.Rows("140:140").ClearFormats()

$ws.Range("A140").Value = "OBI:0000071"
$ws.Range("B140").Value = "quantitative confidence value"
$ws.Range("C140").Value = "A data item which is used to indicate the degree of uncertainty about a measurement."
$ws.Range("D140").Value = "data item"
$ws.Range("J140").Value = "Intervention outcomes and spillover effects"
$ws.Range("P140").Value = "LSR 1"
$ws.Range("Q140").Value = "Intervention outcomes and spillover effects"
$ws.Range("S140").Value = "External"
$ws.Range("V140").Value = "PS"
